$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new headers I1 (I0) and J1 (IF), copying the style from H1 ---
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Fill in the I and J columns for rows 2..85 ---
$iVals = @(9,9,9,9,9,9,9,10,9,9,9,10,9,10,9,10,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,10,9,9,9,9,9,8,9,10,9,9,9,9,8,9,9,10,9,9,9,9,9,9,9,9,10,9,9,9,9,9,9,9,9,9,9,10,9,9,9,8,3,9,5,9,5,5,4,3)
$jVals = @(9,9,9,9,9,9,9,10,9,9,9,10,9,10,9,11,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,10,9,9,9,9,9,9,9,10,9,9,9,9,8,9,9,10,9,9,9,9,9,9,9,9,10,9,9,9,9,9,9,9,9,9,9,10,9,9,9,8,3,9,5,9,5,5,4,3)

for ($idx = 0; $idx -lt $iVals.Count; $idx++) {
  $row = $idx + 2
  $ws.Cells.Item($row, 9).Value = $iVals[$idx]
  $ws.Cells.Item($row, 10).Value = $jVals[$idx]
}
